$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) "Flavio - Conception du trello et modification       Readme"
#    Merge the split "trello"/"Readme" runs (with proofErr spell-check
#    markers) back into a single run, while keeping "Flavio" as its own
#    run. Also drop the paragraph's now-stale <w:pPr><w:rPr> block.
# ---------------------------------------------------------------
$d.Content.Find.Execute(" - Conception du trello et modification       Readme", $true, $false, $false, $false, $false, `
    $true, 1, $false, " - Conception du trello et modification       Readme", 2) | Out-Null

# The replace above merges every run of that paragraph (including the
# leading "Flavio" run) into one. Re-split "Flavio" back into its own
# run by toggling a character property across just that span and back
# to its original value - Word only needs to break the run apart, the
# resulting <w:rPr> ends up identical either side.
$i = 0
$flavioPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Flavio - Conception du trello et modification       Readme") {
        $flavioPara = $p
    }
    $i = $i + 1
}
$pr = $flavioPara.Range
$flavioRange = $d.Range($pr.Start, $pr.Start + 6)
$flavioRange.Bold = 1
$flavioRange.Bold = 0

# Finally strip the leftover paragraph-mark run properties (<w:pPr><w:rPr>)
$flavioPara.Range.ParagraphFormat.Reset()

# ---------------------------------------------------------------
# 2) "06/02/2023" - merge the four split runs "06" + "/0" + "2" + "/2023"
# ---------------------------------------------------------------
$d.Content.Find.Execute("06/02/2023", $true, $false, $false, $false, $false, `
    $true, 1, $false, "06/02/2023", 2) | Out-Null

# ---------------------------------------------------------------
# 3) "...faire le product backlog" - merge the split "product"/"backlog"
#    runs (with proofErr spell-check markers) back into one run.
# ---------------------------------------------------------------
$d.Content.Find.Execute("On a comme objectif terminer les maquettes, le cahier des charges et faire le product backlog", $true, $false, $false, $false, $false, `
    $true, 1, $false, "On a comme objectif terminer les maquettes, le cahier des charges et faire le product backlog", 2) | Out-Null
